$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Replace the old fan-blower manufacturer/part-number with the new motor part
$ws.Range("B7").Value = "STEADY MOTOR"
$ws.Range("C7").Value = "WM7040-24V"

# Match the author's final selection on the sheet
[void]$ws.Range("C8").Select()
